$wb = $excel.ActiveWorkbook
$config = $wb.Worksheets.Item("Configuration")
$config.Range("C15").Value = 300
$config.Range("C16").Value = 300
$config.Range("F15").Value = "Mk2A 300x300"
$config.Range("I15").Value = "Double"
